# Fixes during Regression Testing
# Updates the "Date" (and in a couple of cases "Result") columns on several
# test-log sheets to reflect the latest regression test run timestamps.

$wb = $excel.ActiveWorkbook

# --- PayNowCC : rows 2-7, column B (Date) ---
$ws = $wb.Worksheets.Item("PayNowCC")
$ws.Range("B2").Value = "Tue Nov 18 02:01:06 IST 2025"
$ws.Range("B3").Value = "Tue Nov 18 02:02:03 IST 2025"
$ws.Range("B4").Value = "Tue Nov 18 02:02:57 IST 2025"
$ws.Range("B5").Value = "Tue Nov 18 02:03:44 IST 2025"
$ws.Range("B6").Value = "Tue Nov 18 02:04:36 IST 2025"
$ws.Range("B7").Value = "Tue Nov 18 02:05:30 IST 2025"

# --- NoModifyAmountCC : row 2, column B (Date) ---
$ws = $wb.Worksheets.Item("NoModifyAmountCC")
$ws.Range("B2").Value = "Tue Nov 18 01:17:29 IST 2025"

# --- NoModifyBillingAddressCC : row 2, column B (Date) ---
$ws = $wb.Worksheets.Item("NoModifyBillingAddressCC")
$ws.Range("B2").Value = "Tue Nov 18 01:22:03 IST 2025"

# --- CCDeferredCC : row 2, column B (Date) ---
$ws = $wb.Worksheets.Item("CCDeferredCC")
$ws.Range("B2").Value = "Tue Nov 18 01:11:36 IST 2025"

# --- CMCAutopayCC : row 2, column B (Date) ---
$ws = $wb.Worksheets.Item("CMCAutopayCC")
$ws.Range("B2").Value = "Tue Nov 18 21:46:39 IST 2025"

# --- PayNowCreditCardDCF : row 2, column A (Result) and B (Date) ---
$ws = $wb.Worksheets.Item("PayNowCreditCardDCF")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Nov 18 01:33:10 IST 2025"

# --- PayNowCreditCardSCF : row 2, column A (Result) and B (Date) ---
$ws = $wb.Worksheets.Item("PayNowCreditCardSCF")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Nov 18 19:24:34 IST 2025"

# --- DCFCCVerbiage : row 2, column B (Date) ---
$ws = $wb.Worksheets.Item("DCFCCVerbiage")
$ws.Range("B2").Value = "Tue Nov 18 02:44:58 IST 2025"

# --- SCFCCVerbiage : row 2, column B (Date) ---
$ws = $wb.Worksheets.Item("SCFCCVerbiage")
$ws.Range("B2").Value = "Tue Nov 18 02:48:37 IST 2025"
